# Bug fix: 'Participant type not present in excel import step'
#
# Inserts a new "Participant Type" header column (with a Current/Past/
# Friend-Family list validation) just before the existing "Gender" column,
# and adds the matching Male/Female list validation to the Gender column
# that the import step already expected but which the sheet never
# enforced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing D1:I1 headers one column to the right (E1:J1), -----
# --- working right-to-left so we never clobber a cell before reading it. --
$ws.Range("J1").Value = $ws.Range("I1").Value2
$ws.Range("J1").Font.Name = $ws.Range("I1").Font.Name
$ws.Range("J1").Font.Bold = $ws.Range("I1").Font.Bold
$ws.Range("J1").Font.Size = $ws.Range("I1").Font.Size

$ws.Range("I1").Value = $ws.Range("H1").Value2
$ws.Range("H1").Value = $ws.Range("G1").Value2
$ws.Range("G1").Value = $ws.Range("F1").Value2
$ws.Range("F1").Value = $ws.Range("E1").Value2
$ws.Range("E1").Value = $ws.Range("D1").Value2

# --- New column D becomes "Participant Type" (keeps D1's existing style). -
$ws.Range("D1").Value = "Participant Type"

# --- Data validation drop-downs matching the import step's expectations. --
$ws.Range("D2:D1048576").Validation.Add(3, 1, 1, '"Current Member, Past Member, Friend/Family of Member"')
$ws.Range("E2:E1048576").Validation.Add(3, 1, 1, '"Male, Female"')

# --- Put the view/selection where the author left it. ---------------------
$ws.Range("E1047684").Select()
